# Exercise 4: rename the original data sheet and add two summary tabs
# (YearlyIncome, MonthlyIncome) with header/total rows styled like the
# rest of the workbook (bold+border headers, currency totals).

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet -------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "OrderDetailsData"

# --- Add the two new sheets, in order, right after OrderDetailsData -----------
$wsYearly = $wb.Worksheets.Add($null, $wsData)
$wsYearly.Name = "YearlyIncome"

$wsMonthly = $wb.Worksheets.Add($null, $wsYearly)
$wsMonthly.Name = "MonthlyIncome"

# --- YearlyIncome sheet ---------------------------------------------------------
$wsYearly.Range("A1").Value = "Year"
$wsYearly.Range("B1").Value = 2016
$wsYearly.Range("C1").Value = 2017
$wsYearly.Range("D1").Value = 2018
$wsYearly.Range("A1:D1").Font.Bold = $true
$wsYearly.Range("A1:D1").Borders.LineStyle = 1

$wsYearly.Range("A2").Value = "Total "
$wsYearly.Range("A2:D2").Borders.LineStyle = 1
$wsYearly.Range("B2:D2").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$wsYearly.Range("B2").Formula = "=SUMIF(OrderDetailsData!D:D,B1,OrderDetailsData!H:H)"
$wsYearly.Range("C2").Formula = "=SUMIF(OrderDetailsData!D:D,C1,OrderDetailsData!H:H)"
$wsYearly.Range("D2").Formula = "=SUMIF(OrderDetailsData!D:D,D1,OrderDetailsData!H:H)"

$wsYearly.Columns.Item(1).ColumnWidth = 10
$wsYearly.Columns.Item(2).ColumnWidth = 13.570312
$wsYearly.Columns.Item(3).ColumnWidth = 17.425781
$wsYearly.Columns.Item(4).ColumnWidth = 17

# --- MonthlyIncome sheet ---------------------------------------------------------
$wsMonthly.Range("A1").Value = "Month"
$wsMonthly.Range("B1").Value = 2016
$wsMonthly.Range("C1").Value = 2017
$wsMonthly.Range("D1").Value = 2018
$wsMonthly.Range("A1:D1").Font.Bold = $true
$wsMonthly.Range("A1:D1").Borders.LineStyle = 1

for ($m = 1; $m -le 12; $m++) {
    $r = $m + 1
    $wsMonthly.Cells.Item($r, 1).Value = $m
    $wsMonthly.Range($wsMonthly.Cells.Item($r, 1), $wsMonthly.Cells.Item($r, 4)).Borders.LineStyle = 1
    $rowRange = $wsMonthly.Range($wsMonthly.Cells.Item($r, 2), $wsMonthly.Cells.Item($r, 4))
    $rowRange.NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
    $wsMonthly.Cells.Item($r, 2).Formula = "=SUMIFS(OrderDetailsData!`$H:`$H,OrderDetailsData!`$D:`$D,B`$1,OrderDetailsData!`$L:`$L,""<=""&DATE(B`$1,$m+1,0),OrderDetailsData!`$L:`$L,"">=""&DATE(B`$1,$m,1))"
    $wsMonthly.Cells.Item($r, 3).Formula = "=SUMIFS(OrderDetailsData!`$H:`$H,OrderDetailsData!`$D:`$D,C`$1,OrderDetailsData!`$L:`$L,""<=""&DATE(C`$1,$m+1,0),OrderDetailsData!`$L:`$L,"">=""&DATE(C`$1,$m,1))"
    $wsMonthly.Cells.Item($r, 4).Formula = "=SUMIFS(OrderDetailsData!`$H:`$H,OrderDetailsData!`$D:`$D,D`$1,OrderDetailsData!`$L:`$L,""<=""&DATE(D`$1,$m+1,0),OrderDetailsData!`$L:`$L,"">=""&DATE(D`$1,$m,1))"
}

$wsMonthly.Range("A14").Value = "Total"
$wsMonthly.Range("A14:D14").Borders.LineStyle = 1
$wsMonthly.Range("A1:A14").Font.Bold = $false
$wsMonthly.Range("A14").Font.Bold = $true
$wsMonthly.Range("B14:D14").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$wsMonthly.Range("B14").Formula = "=SUM(B2:B13)"
$wsMonthly.Range("C14").Formula = "=SUM(C2:C13)"
$wsMonthly.Range("D14").Formula = "=SUM(D2:D13)"

$wsMonthly.Columns.Item(1).ColumnWidth = 10
$wsMonthly.Columns.Item(2).ColumnWidth = 13.570312
$wsMonthly.Columns.Item(3).ColumnWidth = 17.425781
$wsMonthly.Columns.Item(4).ColumnWidth = 17

# --- Restore the active sheet / selection on the data tab ----------------------
$wsData.Activate()
$wsData.Range("A1").Select()

Write-Output "done"
